$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3616.25
$ws.Range("I62").Value = 6300
$ws.Range("J62").Value = 2006
$ws.Range("K62").Value = 6300
$ws.Range("L62").Value = 2006
$ws.Range("M62").Value = -5676
$ws.Range("N62").Value = -3254

$ws.Range("H65").Value = 3616.25
$ws.Range("I65").Value = 6300
$ws.Range("J65").Value = 2006
$ws.Range("K65").Value = 31500
$ws.Range("L65").Value = 10030
$ws.Range("M65").Value = -28380
$ws.Range("N65").Value = -16270

$ws.Range("H98").Value = 23385.4
$ws.Range("I98").Value = 786
$ws.Range("J98").Value = 53517.934
$ws.Range("K98").Value = 786
$ws.Range("L98").Value = 53517.934
$ws.Range("M98").Value = 712
$ws.Range("N98").Value = -56513.934

$ws.Range("H117").Value = 44425.5
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 44425.5
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 44425.5
$ws.Range("N117").Value = -53603.5

$ws.Range("H122").Value = 23385.4
$ws.Range("I122").Value = 786
$ws.Range("J122").Value = 53517.934
$ws.Range("K122").Value = 2358
$ws.Range("L122").Value = 160553.802
$ws.Range("M122").Value = 92
$ws.Range("N122").Value = -165453.802

$ws.Range("H129").Value = 1524.5476
$ws.Range("I129").Value = 3073.5
$ws.Range("J129").Value = 1361.5
$ws.Range("K129").Value = 9220.5
$ws.Range("L129").Value = 4084.5
$ws.Range("M129").Value = -4220.5
$ws.Range("N129").Value = -14084.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22367.756
$ws.Range("I32").Value = 20739.096
$ws.Range("J32").Value = 45169
$ws.Range("K32").Value = 20739.096
$ws.Range("L32").Value = 45169
$ws.Range("M32").Value = -20452.096
$ws.Range("N32").Value = -45743

$ws.Range("H110").Value = 1248.2222
$ws.Range("I110").Value = 1029.25
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 1029.25
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = 1015.75
$ws.Range("N110").Value = -7090

$ws.Range("H122").Value = 1809.7894
$ws.Range("I122").Value = 1809.7894
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5429.3682
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2979.3682

$ws.Range("H131").Value = 46670.75
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 46670.75
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 46670.75
$ws.Range("N131").Value = -56750.75

$ws.Range("H132").Value = 15627314
$ws.Range("I132").Value = 26317290
$ws.Range("J132").Value = 3503.3845
$ws.Range("K132").Value = 78951870
$ws.Range("L132").Value = 10510.1535
$ws.Range("M132").Value = -78949340
$ws.Range("N132").Value = -15570.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1640.2812
$ws.Range("I99").Value = 1543.2963
$ws.Range("J99").Value = 2164
$ws.Range("K99").Value = 1543.2963
$ws.Range("L99").Value = 2164
$ws.Range("M99").Value = -45.29629999999997
$ws.Range("N99").Value = -5160

$ws.Range("H108").Value = 48676
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 48676
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 48676
$ws.Range("N108").Value = -56356

$ws.Range("H130").Value = 48731.8
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 48731.8
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 48731.8
$ws.Range("N130").Value = -58771.8

$ws.Range("H141").Value = 39893
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 39893
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 39893
$ws.Range("N141").Value = -50253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 35512.8
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 35512.8
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 35512.8
$ws.Range("N28").Value = -36002.8

$ws.Range("H31").Value = 4909.397
$ws.Range("I31").Value = 1751.4762
$ws.Range("J31").Value = 6488.357
$ws.Range("K31").Value = 1751.4762
$ws.Range("L31").Value = 6488.357
$ws.Range("M31").Value = -1456.4762
$ws.Range("N31").Value = -7078.357

$ws.Range("H34").Value = 4909.397
$ws.Range("I34").Value = 1751.4762
$ws.Range("J34").Value = 6488.357
$ws.Range("K34").Value = 1751.4762
$ws.Range("L34").Value = 6488.357
$ws.Range("M34").Value = -1549.4762
$ws.Range("N34").Value = -6892.357

$ws.Range("H43").Value = 147056.83
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 147056.83
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 147056.83
$ws.Range("N43").Value = -147424.83

$ws.Range("H95").Value = 69833.336
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 69833.336
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 69833.336
$ws.Range("N95").Value = -75325.336

$ws.Range("H101").Value = 147056.83
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 147056.83
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 147056.83
$ws.Range("N101").Value = -153546.83

$ws.Range("H132").Value = 63536.086
$ws.Range("I132").Value = 2224
$ws.Range("J132").Value = 158910.44
$ws.Range("K132").Value = 6672
$ws.Range("L132").Value = 476731.32
$ws.Range("M132").Value = -4142
$ws.Range("N132").Value = -481791.32

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1153.75
$ws.Range("I92").Value = 1037.8334
$ws.Range("J92").Value = 1501.5
$ws.Range("K92").Value = 3113.5002
$ws.Range("L92").Value = 4504.5
$ws.Range("M92").Value = -1865.5002
$ws.Range("N92").Value = -7000.5

$ws.Range("H113").Value = 3073.775
$ws.Range("I113").Value = 3995.6206
$ws.Range("J113").Value = 643.4545000000001
$ws.Range("K113").Value = 11986.8618
$ws.Range("L113").Value = 1930.3635
$ws.Range("M113").Value = -9816.861800000001
$ws.Range("N113").Value = -6270.3635

$ws.Range("H118").Value = 3717.6875
$ws.Range("I118").Value = 2286.3333
$ws.Range("J118").Value = 4048
$ws.Range("K118").Value = 6858.999899999999
$ws.Range("L118").Value = 12144
$ws.Range("M118").Value = -5615.999899999999
$ws.Range("N118").Value = -14630

$ws.Range("H131").Value = 61249.855
$ws.Range("I131").Value = 20641.4
$ws.Range("J131").Value = 68017.92999999999
$ws.Range("K131").Value = 61924.2
$ws.Range("L131").Value = 204053.79
$ws.Range("M131").Value = -56884.2
$ws.Range("N131").Value = -214133.79

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 34496
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 34496
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 34496
$ws.Range("N110").Value = -42676

$ws.Range("H130").Value = 53000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 53000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 53000
$ws.Range("N130").Value = -63040

$ws.Range("H133").Value = 26783.691
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 26783.691
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 26783.691
$ws.Range("N133").Value = -36903.691

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2290.5715
$ws.Range("I40").Value = 2262.6428
$ws.Range("J40").Value = 2346.4285
$ws.Range("K40").Value = 2262.6428
$ws.Range("L40").Value = 2346.4285
$ws.Range("M40").Value = -2126.6428
$ws.Range("N40").Value = -2618.4285

$ws.Range("H108").Value = 40414.668
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 40414.668
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 40414.668
$ws.Range("N108").Value = -48094.668

$ws.Range("H136").Value = 3178.7368
$ws.Range("I136").Value = 2484.3845
$ws.Range("J136").Value = 4683.1665
$ws.Range("K136").Value = 7453.1535
$ws.Range("L136").Value = 14049.4995
$ws.Range("M136").Value = -4903.1535
$ws.Range("N136").Value = -19149.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 44989
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 44989
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 44989
$ws.Range("N16").Value = -45573

$ws.Range("H39").Value = 23000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 23000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 23000
$ws.Range("N39").Value = -23826

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H113").Value = 783.3333
$ws.Range("I113").Value = 740
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2220
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -50
$ws.Range("N113").Value = -7340

$ws.Range("H119").Value = 49690
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 49690
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 49690
$ws.Range("N119").Value = -59366

$ws.Range("H122").Value = 1361478.2
$ws.Range("I122").Value = 2041582
$ws.Range("J122").Value = 1270.7142
$ws.Range("K122").Value = 6124746
$ws.Range("L122").Value = 3812.1426
$ws.Range("M122").Value = -6122296
$ws.Range("N122").Value = -8712.142599999999

$ws.Range("H132").Value = 3697.4285
$ws.Range("I132").Value = 3196.6
$ws.Range("J132").Value = 4949.5
$ws.Range("K132").Value = 9589.799999999999
$ws.Range("L132").Value = 14848.5
$ws.Range("M132").Value = -7059.799999999999
$ws.Range("N132").Value = -19908.5
